# Update the cryptocurrency price/volume table to reflect refreshed data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.936.26'
$ws.Range("E2").Value = '  +2.63%  '
$ws.Range("D3").Value = '3.568.83'
$ws.Range("E3").Value = '  +1.52%  '
$ws.Range("E4").Value = '  -0.06%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '583.09'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +2.07%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '186.18'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +2.52%  '
$ws.Range("D7").Value = '3.556.77'
$ws.Range("E7").Value = '  +1.36%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.621'
$c.Style = "Normal"
$ws.Range("E8").Value = '  +1.37%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E9").Value = '  +0.09%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.216'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +13.96%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.652'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +2.57%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '54.78'
$c.Style = "Normal"
$ws.Range("E12").Value = '  +2.18%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '0.0000318'
$c.Style = "Normal"
$ws.Range("E13").Value = '  +5.52%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '9.57'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +1.39%  '
$ws.Range("D15").Value = '4.046.77'
$ws.Range("E15").Value = '  -1.24%  '
$ws.Range("D16").Value = '71.021.26'
$ws.Range("E16").Value = '  +2.82%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.602.81'
$ws.Range("E17").Value = '  +2.34%  '
$ws.Range("B18").Value = 'Chainlink'
$ws.Range("C18").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '19.32'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +0.49%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '12.40'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -0.25%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '563.73'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +6.12%  '
$ws.Range("E21").Value = '  +0.69%  '
$ws.Range("E22").Value = '  -1.49%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '17.65'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -11.10%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '5.02'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +1.38%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '4.57'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +4.89%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '94.57'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +0.74%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '11.35'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +3.67%  '
$ws.Range("E28").Value = '  +1.99%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '9.18'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +1.23%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '32.54'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +3.33%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '7.32'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +0.77%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '12.34'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -1.63%  '
$ws.Range("E33").Value = '  +2.48%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '63.74'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -0.80%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '3.42'
$c.Style = "Normal"
$ws.Range("E35").Value = '  +10.56%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '553.98'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -2.72%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '0.420'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +5.87%  '
$ws.Range("D38").Value = '0.0₃0802'
$ws.Range("E38").Value = '  +5.48%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '37.70'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -0.63%  '
$ws.Range("E40").Value = '  -0.06%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '3.31'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +8.69%  '
$ws.Range("D42").Value = '3.550.85'
$ws.Range("E42").Value = '  +12.24%  '
$ws.Range("B43").Value = 'Kaspa'
$ws.Range("C43").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.137'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +3.08%  '
$ws.Range("B44").Value = 'Stacks'
$ws.Range("C44").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '3.43'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +2.48%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.0449'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +1.77%  '
$ws.Range("B46").Value = 'ApeXProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '3.48'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -1.02%  '
$ws.Range("B47").Value = 'ThetaToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '2.95'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -0.41%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '9.40'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +2.02%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '0.137'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +3.06%  '
$ws.Range("E50").Value = '  +8.80%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -0.03%  '
